$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free direct cell updates reproducing the diff row by row.
# For D-column values that would otherwise be auto-parsed as numbers by
# Excel (plain decimals with a single "." separator), we briefly force a
# text number format, set the literal text, then clear the format again so
# the cell keeps its original (un-styled) appearance but literal text value.

$ws.Range("D2").Value = '28.466.24'
$ws.Range("E2").Value = '  +1.70%  '

$ws.Range("D3").Value = '1.826.06'
$ws.Range("E3").Value = '  +1.56%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.38'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.42%  '

$ws.Range("E6").Value = '  -0.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5140'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -5.38%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3939'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +3.70%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07699'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.94'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.112'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.87%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.97'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.21%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.278'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.08%  '

$ws.Range("E14").Value = '  +0.00%  '

$ws.Range("E15").Value = '  +2.26%  '

$ws.Range("D16").Value = '1.825.70'
$ws.Range("E16").Value = '  +1.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.61'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +5.21%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001080'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.57%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06633'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.76%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.67'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.59%  '

$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("E22").Value = '  +2.06%  '

$ws.Range("D23").Value = '28.493.21'
$ws.Range("E23").Value = '  +1.66%  '

$ws.Range("E24").Value = '  -0.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.242'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +7.25%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.03'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.40%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.432'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.90%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.59'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.22%  '

$ws.Range("D29").Value = '2.033.52'
$ws.Range("E29").Value = '  +1.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.95'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.128'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1098'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.68%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.650'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.650'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.57%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07167'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.36%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2231'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.986'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +6.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02325'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.56%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.162'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.67%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6235'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +1.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.25'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.10%  '

$ws.Range("E42").Value = '  +1.65%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.394'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.72%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.45'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5892'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.90%  '

$ws.Range("E47").Value = '  +0.55%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.38'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.976'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.182'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06940'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.90%  '
